$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Modules" header row: swap the Enseignant/Nombre d'heures columns for
# Composants / Chef  Module (C1 <-> D1 semantics change).
$ws.Range("C1").Value = "Chef  Module"
$ws.Range("D1").Value = "Composants"

# Widen the two columns so the new (longer) headers are readable.
$ws.Columns.Item(3).ColumnWidth = 34.1666666666667
$ws.Columns.Item(4).ColumnWidth = 23.6666666666667

# Move the active selection.
[void]$ws.Range("E8").Select()
